# BOM.xlsx edit: "Changed horizontal button position" / resistor value update
# and addition of a new tool (electronics tweezers) to the BOM tools list.
#
# The workbook's ActiveSheet (the BOM data sheet, tabSelected) already
# resolves correctly via $wb.ActiveSheet - use that instead of chasing
# sheet names/activation, since the physical sheet1.xml / sheet2.xml file
# names don't line up with the worksheet tab names in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the tools list with a new row for "electronics tweezers" ---
# Row 28 ("loupe tool") is the last existing tools row; copy its cell
# formatting down into the new row 29 before filling in the new row's
# values.
$ws.Range("A28").Copy()
$ws.Range("A29").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C28").Copy()
$ws.Range("C29").PasteSpecial(-4122)   # xlPasteFormats

# Fill the new row content first so "electronics tweezers" is registered
# in the shared-string table ahead of the resistor value change below.
$ws.Range("C29").Value = "electronics tweezers"
$ws.Range("B29").Value = 1

# The "loupe tool" row (28) also picks up a quantity of 1.
$ws.Range("B28").Value = 1

# --- Resistor value correction on row 8: 51R -> 56R ---
$ws.Range("C8").Value = "56R"

# Leave the selection where the editor ended up after making the change.
$ws.Range("C9").Select()
